# Fruta / hortaliza, semanal
# Insert a new weekly price row at row 170 (pushing the existing rows 170-180
# down to 171-181) on the "Piña" price log worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 170, shifting rows 170-180
# down to become rows 171-181 (this also grows the used range to A1:T181).
$ws.Rows("170:170").Insert()

# Populate the newly inserted row 170 with this week's data.
$ws.Range("A170").Value = 11
$ws.Range("B170").Value = "Vega Monumental Concepción"
$ws.Range("C170").Value = "Bíobío"
$ws.Range("D170").Value = (Get-Date -Year 2022 -Month 7 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E170").Value = 8
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100108
$ws.Range("H170").Value = "Tropicales y subtropicales"
$ws.Range("I170").Value = 100108005
$ws.Range("J170").Value = "Piña"
$ws.Range("K170").Value = "Caramelo"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 100
$ws.Range("N170").Value = 19000
$ws.Range("O170").Value = 20000
$ws.Range("P170").Value = 19500
$ws.Range("Q170").Value = '$/caja 14 unidades'
$ws.Range("R170").Value = "Ecuador"
$ws.Range("S170").Value = 1393
$ws.Range("T170").Value = 14
